$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.04882852104406
$ws.Range("D2").Value = 1.058826362515135
$ws.Range("E2").Value = 1.05629583397662
$ws.Range("F2").Value = 1.067509638405109
$ws.Range("I2").Value = 1.036406751670986
$ws.Range("J2").Value = 1.053870200571587
$ws.Range("K2").Value = 1.061557536056161
$ws.Range("L2").Value = 1.059033935497222
$ws.Range("M2").Value = 1.070217309409141
$ws.Range("N2").Value = 1.055366817137001

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.049891122660022
$ws.Range("D3").Value = 1.059832868492836
$ws.Range("E3").Value = 1.057255794729825
$ws.Range("F3").Value = 1.068587750771924
$ws.Range("I3").Value = 1.036539604502645
$ws.Range("J3").Value = 1.054581282238466
$ws.Range("K3").Value = 1.062377664003399
$ws.Range("L3").Value = 1.059807134467283
$ws.Range("M3").Value = 1.071110568008695
$ws.Range("N3").Value = 1.056078908621408

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.050578959081715
$ws.Range("D4").Value = 1.060484730195321
$ws.Range("E4").Value = 1.057877569342166
$ws.Range("F4").Value = 1.069286178914241
$ws.Range("I4").Value = 1.036624076224236
$ws.Range("J4").Value = 1.055041071068849
$ws.Range("K4").Value = 1.062908312681881
$ws.Range("L4").Value = 1.060307438448018
$ws.Range("M4").Value = 1.071688781656671
$ws.Range("N4").Value = 1.056539350404669

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.050868188130886
$ws.Range("D5").Value = 1.060758912318741
$ws.Range("E5").Value = 1.058139110171916
$ws.Range("F5").Value = 1.069579993215742
$ws.Range("I5").Value = 1.036659230369237
$ws.Range("J5").Value = 1.055234287177611
$ws.Range("K5").Value = 1.063131389952974
$ws.Range("L5").Value = 1.060517763954234
$ws.Range("M5").Value = 1.071931913767783
$ws.Range("N5").Value = 1.056732840902468

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.050916754629007
$ws.Range("D6").Value = 1.060804956864716
$ws.Range("E6").Value = 1.058183032616775
$ws.Range("F6").Value = 1.069629337354825
$ws.Range("I6").Value = 1.036665111915
$ws.Range("J6").Value = 1.055266724359797
$ws.Range("K6").Value = 1.063168845157865
$ws.Range("L6").Value = 1.060553078383502
$ws.Range("M6").Value = 1.071972739700076
$ws.Range("N6").Value = 1.056765324149173

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.050582823528544
$ws.Range("D7").Value = 1.06048839328302
$ws.Range("E7").Value = 1.057881063487586
$ws.Range("F7").Value = 1.069290104109372
$ws.Range("I7").Value = 1.03662454736186
$ws.Range("J7").Value = 1.055043653143397
$ws.Range("K7").Value = 1.062911293481136
$ws.Range("L7").Value = 1.060310248837344
$ws.Range("M7").Value = 1.071692030199732
$ws.Range("N7").Value = 1.056541936146059

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.049187577879113
$ws.Range("D8").Value = 1.059166394430653
$ws.Range("E8").Value = 1.056620129311719
$ws.Range("F8").Value = 1.067873821949647
$ws.Range("I8").Value = 1.03645195881928
$ws.Range("J8").Value = 1.054110581704348
$ws.Range("K8").Value = 1.061834707825495
$ws.Range("L8").Value = 1.059295242950705
$ws.Range("M8").Value = 1.070519145378617
$ws.Range("N8").Value = 1.055607539638542

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.046730989929834
$ws.Range("D9").Value = 1.056841370129058
$ws.Range("E9").Value = 1.054402946771389
$ws.Range("F9").Value = 1.065384436394245
$ws.Range("I9").Value = 1.036136417069307
$ws.Range("J9").Value = 1.052463894238444
$ws.Range("K9").Value = 1.059937426660132
$ws.Range("L9").Value = 1.057506637725867
$ws.Range("M9").Value = 1.068454046524498
$ws.Range("N9").Value = 1.053958513687558

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.045094622896023
$ws.Range("D10").Value = 1.055294416718524
$ws.Range("E10").Value = 1.052928051092529
$ws.Range("F10").Value = 1.063729102356286
$ws.Range("I10").Value = 1.035918398876922
$ws.Range("J10").Value = 1.051364449320434
$ws.Range("K10").Value = 1.058672465996577
$ws.Range("L10").Value = 1.056314238365273
$ws.Range("M10").Value = 1.067078468195239
$ws.Range("N10").Value = 1.052857507431657

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.044386380138244
$ws.Range("D11").Value = 1.054625301487968
$ws.Range("E11").Value = 1.052290177107781
$ws.Range("F11").Value = 1.063013339776131
$ws.Range("I11").Value = 1.035822182500949
$ws.Range("J11").Value = 1.050887990113396
$ws.Range("K11").Value = 1.058124703999186
$ws.Range("L11").Value = 1.055797922523215
$ws.Range("M11").Value = 1.066483106752649
$ws.Range("N11").Value = 1.052380371597894

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.044123354455663
$ws.Range("D12").Value = 1.054376871763306
$ws.Range("E12").Value = 1.052053357815836
$ws.Range("F12").Value = 1.062747625505037
$ws.Range("I12").Value = 1.035786171346674
$ws.Range("J12").Value = 1.050710953314083
$ws.Range("K12").Value = 1.057921237286202
$ws.Range("L12").Value = 1.055606140324486
$ws.Range("M12").Value = 1.066262004381355
$ws.Range("N12").Value = 1.05220308338602

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.04417977222268
$ws.Range("D13").Value = 1.054430155822305
$ws.Range("E13").Value = 1.052104151123066
$ws.Range("F13").Value = 1.062804615241815
$ws.Range("I13").Value = 1.035793908170422
$ws.Range("J13").Value = 1.050748930956994
$ws.Range("K13").Value = 1.05796488174054
$ws.Range("L13").Value = 1.055647278229129
$ws.Range("M13").Value = 1.066309429712452
$ws.Range("N13").Value = 1.05224111496154

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.044364637390764
$ws.Range("D14").Value = 1.054604763973953
$ws.Range("E14").Value = 1.052270599192091
$ws.Range("F14").Value = 1.06299137265228
$ws.Range("I14").Value = 1.035819211356986
$ws.Range("J14").Value = 1.050873357390723
$ws.Range("K14").Value = 1.058107885447011
$ws.Range("L14").Value = 1.055782069734728
$ws.Range("M14").Value = 1.066464829515533
$ws.Range("N14").Value = 1.052365718095076

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.044478545175061
$ws.Range("D15").Value = 1.054712360336667
$ws.Range("E15").Value = 1.052373168661356
$ws.Range("F15").Value = 1.063106460185812
$ws.Range("I15").Value = 1.035834765436851
$ws.Range("J15").Value = 1.050950012862467
$ws.Range("K15").Value = 1.058195994294322
$ws.Range("L15").Value = 1.055865119305996
$ws.Range("M15").Value = 1.066560581956522
$ws.Range("N15").Value = 1.052442482426383

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.045141633261985
$ws.Range("D16").Value = 1.055338839020421
$ws.Range("E16").Value = 1.052970400859578
$ws.Range("F16").Value = 1.063776626431061
$ws.Range("I16").Value = 1.03592474626602
$ws.Range("J16").Value = 1.051396062109348
$ws.Range("K16").Value = 1.058708818644263
$ws.Range("L16").Value = 1.056348504597237
$ws.Range("M16").Value = 1.067117986156476
$ws.Range("N16").Value = 1.052889165114357

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.045557655148536
$ws.Range("D17").Value = 1.055732007162434
$ws.Range("E17").Value = 1.053345234150022
$ws.Range("F17").Value = 1.064197274146244
$ws.Range("I17").Value = 1.035980703655377
$ws.Range("J17").Value = 1.051675752201061
$ws.Range("K17").Value = 1.059030493174503
$ws.Range("L17").Value = 1.056651719921657
$ws.Range("M17").Value = 1.067467704437969
$ws.Range("N17").Value = 1.053169252398081

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.045800344073138
$ws.Range("D18").Value = 1.055961405489156
$ws.Range("E18").Value = 1.053563941900985
$ws.Range("F18").Value = 1.064442728071598
$ws.Range("I18").Value = 1.036013167700836
$ws.Range("J18").Value = 1.051838852827678
$ws.Range("K18").Value = 1.059218117984794
$ws.Range("L18").Value = 1.056828580358418
$ws.Range("M18").Value = 1.067671715516483
$ws.Range("N18").Value = 1.053332584646296

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.045883099883932
$ws.Range("D19").Value = 1.056039636251055
$ws.Range("E19").Value = 1.053638528100891
$ws.Range("F19").Value = 1.064526437997303
$ws.Range("I19").Value = 1.036024207417701
$ws.Range("J19").Value = 1.051894459493508
$ws.Range("K19").Value = 1.059282092761908
$ws.Range("L19").Value = 1.056888885197735
$ws.Range("M19").Value = 1.067741282476001
$ws.Range("N19").Value = 1.05338827027997

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.045513016793887
$ws.Range("D20").Value = 1.055689816698006
$ws.Range("E20").Value = 1.053305010444707
$ws.Range("F20").Value = 1.064152132567573
$ws.Range("I20").Value = 1.035974718051214
$ws.Range("J20").Value = 1.051645748016179
$ws.Range("K20").Value = 1.058995980797843
$ws.Range("L20").Value = 1.056619187777297
$ws.Range("M20").Value = 1.067430180250846
$ws.Range("N20").Value = 1.053139205603815

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.044310197903967
$ws.Range("D21").Value = 1.054553343200154
$ws.Range("E21").Value = 1.052221581191882
$ws.Range("F21").Value = 1.062936373055634
$ws.Range("I21").Value = 1.035811767708448
$ws.Range("J21").Value = 1.050836718523421
$ws.Range("K21").Value = 1.058065774514104
$ws.Range("L21").Value = 1.055742376979457
$ws.Range("M21").Value = 1.066419067001708
$ws.Range("N21").Value = 1.05232902719638

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.04355421103275
$ws.Range("D22").Value = 1.05383943074145
$ws.Range("E22").Value = 1.051541054751427
$ws.Range("F22").Value = 1.062172855374525
$ws.Range("I22").Value = 1.035707740166954
$ws.Range("J22").Value = 1.050327710475379
$ws.Range("K22").Value = 1.057480896718902
$ws.Range("L22").Value = 1.055191094263007
$ws.Range("M22").Value = 1.065783579541666
$ws.Range("N22").Value = 1.051819296298526

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.043954948124807
$ws.Range("D23").Value = 1.054217829050929
$ws.Range("E23").Value = 1.051901751161271
$ws.Range("F23").Value = 1.062577527136331
$ws.Range("I23").Value = 1.035763036255559
$ws.Range("J23").Value = 1.05059757731315
$ws.Range("K23").Value = 1.057790953303703
$ws.Range("L23").Value = 1.055483339256619
$ws.Range("M23").Value = 1.066120440640959
$ws.Range("N23").Value = 1.052089546378161

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.045533186850834
$ws.Range("D24").Value = 1.055708880536461
$ws.Range("E24").Value = 1.053323185577553
$ws.Range("F24").Value = 1.064172529803384
$ws.Range("I24").Value = 1.035977423228663
$ws.Range("J24").Value = 1.051659305732768
$ws.Range("K24").Value = 1.059011575462952
$ws.Range("L24").Value = 1.056633887652739
$ws.Range("M24").Value = 1.067447135734403
$ws.Range("N24").Value = 1.053152782573917

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.047365837753501
$ws.Range("D25").Value = 1.057441907209152
$ws.Range("E25").Value = 1.054975575634982
$ws.Range("F25").Value = 1.06602725429319
$ws.Range("I25").Value = 1.036219343261541
$ws.Range("J25").Value = 1.052889895663125
$ws.Range("K25").Value = 1.060427940130624
$ws.Range("L25").Value = 1.057969036357195
$ws.Range("M25").Value = 1.070217309409141
$ws.Range("N25").Value = 1.054385120083125
